$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Rspo2"
$ws.Range("C2").Value = "Lgr5"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.621828666666667
$ws.Range("H2").Value = 4.865486
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.1978116666666667
$ws.Range("N2").Value = 0.593435
$ws.Range("O2").Value = 0.0447866041160341
$ws.Range("P2").Value = 0.04478660411603409
$ws.Range("Q2").Value = 0.3208166316011111
$ws.Range("R2").Value = 2.88734968441
$ws.Range("S2").Value = 0.0447866041160341
$ws.Range("T2").Value = 0.04478660411603409
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Rspo2"
$ws.Range("C3").Value = "Lgr5"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.621828666666667
$ws.Range("H3").Value = 4.865486
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.940645333333333
$ws.Range("N3").Value = 5.821936
$ws.Range("O3").Value = 0.4393821443306968
$ws.Range("P3").Value = 0.4393821443306967
$ws.Range("Q3").Value = 3.147394233432889
$ws.Range("R3").Value = 28.326548100896
$ws.Range("S3").Value = 0.4393821443306968
$ws.Range("T3").Value = 0.4393821443306967
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Rspo2"
$ws.Range("C4").Value = "Lgr5"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.621828666666667
$ws.Range("H4").Value = 4.865486
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.046181
$ws.Range("N4").Value = 0.138543
$ws.Range("O4").Value = 0.0104558553068958
$ws.Range("P4").Value = 0.0104558553068958
$ws.Range("Q4").Value = 0.07489766965533333
$ws.Range("R4").Value = 0.6740790268979999
$ws.Range("S4").Value = 0.0104558553068958
$ws.Range("T4").Value = 0.0104558553068958
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Rspo2"
$ws.Range("C5").Value = "Lgr5"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.621828666666667
$ws.Range("H5").Value = 4.865486
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 2.232121666666667
$ws.Range("N5").Value = 6.696365
$ws.Range("O5").Value = 0.5053753962463735
$ws.Range("P5").Value = 0.5053753962463734
$ws.Range("Q5").Value = 3.620118906487777
$ws.Range("R5").Value = 32.58107015839
$ws.Range("S5").Value = 0.5053753962463735
$ws.Range("T5").Value = 0.5053753962463734
